# Update the dSF (column F) values for the affected rows, per the
# "repull data, push all data, mean calculation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 6
    3  = -7
    5  = -5
    6  = 3
    7  = 3
    14 = 2
    15 = 0
    17 = 2
    23 = 0
    24 = 0
    25 = 0
    43 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
